# Upload excel files with prices
# This script applies the scraped-data refresh to the bread_coop workbook:
#   - updates the timestamp column (O) for every data row to the new crawl time
#   - bumps ratingAmount (D) / ratingValue (E) for products whose ratings changed
#   - updates the productAriaLabel (M) text for the out-of-stock item

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Refresh the crawl timestamp (column O) for all data rows (2..395) ---
$oldTimestamp = "2023-01-26 06:49:16"
$newTimestamp = "2023-01-26 12:57:39"

$lastRow = $ws.Cells.Item($ws.Rows.Count, 15).End(-4162).Row
if ($lastRow -lt 395) { $lastRow = 395 }

for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 15).Value = $newTimestamp
}

# --- Rating amount / rating value updates (columns D and E) ---
$ratingChanges = @(
    @{ Row = 9;   D = 30 },
    @{ Row = 23;  D = 23 },
    @{ Row = 30;  D = 7;  E = 4.5 },
    @{ Row = 31;  D = 14 },
    @{ Row = 39;  D = 17 },
    @{ Row = 47;  D = 11 },
    @{ Row = 48;  D = 24 },
    @{ Row = 49;  D = 17 },
    @{ Row = 120; D = 8;  E = 5 },
    @{ Row = 125; D = 43 },
    @{ Row = 165; D = 9 },
    @{ Row = 179; D = 16 },
    @{ Row = 298; D = 4;  E = 5 },
    @{ Row = 299; D = 3;  E = 4.5 }
)

foreach ($change in $ratingChanges) {
    $ws.Cells.Item($change.Row, 4).Value = $change.D
    if ($change.ContainsKey("E")) {
        $ws.Cells.Item($change.Row, 5).Value = $change.E
    }
}

# --- Product aria label text update for row 96 (now "Online kein Bestand") ---
$ws.Range("M96").Value = "Betty Bossi Naturaplan Bio Blätterteig ausgewallt - Online kein Bestand 2.40 Schweizer Franken"
